$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 144, shifting existing rows 144:157 down to 145:158
$ws.Rows.Item(144).Insert()

# Populate the newly-inserted row 144 with the new record
$ws.Range("A144").Value = 4
$ws.Range("B144").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C144").Value = "Los Lagos"
$ws.Range("D144").Value = 44449
$ws.Range("E144").Value = 10
$ws.Range("F144").Value = 100112037
$ws.Range("G144").Value = "Cebollín"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 150
$ws.Range("K144").Value = 6500
$ws.Range("L144").Value = 6500
$ws.Range("M144").Value = 6500
$ws.Range("N144").Value = "`$/paquete 36 unidades"
$ws.Range("O144").Value = "Región Metropolitana"
$ws.Range("P144").Value = 181
$ws.Range("Q144").Value = 36
$ws.Range("R144").Value = "Hortaliza"
